# ---------------------------------------------------------------------------
# Final cleanup: manifest, duplicate files, and Excel names
#
#  1. Disambiguate the "Est. Monthly Churn" description on the Inputs sheet
#     into a Basic-specific and a Pro-specific label.
#  2. Add two new model parameters (Activation Rate, Conversion Rate) as
#     rows 9 and 10 of the Inputs sheet, including their Key/Description
#     columns.
#  3. Register those two new parameters as workbook-level defined names
#     (activation_rate, conversion_rate) pointing at Inputs!$B$9/$B$10.
#  4. Add a new "__names__" worksheet at the end of the workbook that
#     documents every defined name and what it refers to.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsInputs = $wb.Worksheets.Item("Inputs")

# --- 1) Disambiguate churn descriptions (Basic vs Pro) ---------------------
$wsInputs.Range("D6").Value = "Est. Monthly Churn B"
$wsInputs.Range("D7").Value = "Est. Monthly Churn P"

# --- 2) Add new parameter rows 9 and 10 -------------------------------------
$wsInputs.Range("A9").Value = "Activation Rate"
$wsInputs.Range("B9").Value = 0.3
$wsInputs.Range("C9").Value = "activation_rate"
$wsInputs.Range("D9").Value = "Est. Activation %"

$wsInputs.Range("A10").Value = "Conversion Rate"
$wsInputs.Range("B10").Value = 0.05
$wsInputs.Range("C10").Value = "conversion_rate"
$wsInputs.Range("D10").Value = "Est. Trial-to-Paid %"

# --- 3) Register the two new workbook-level defined names -------------------
$wb.Names.Add('activation_rate', '=Inputs!$B$9')
$wb.Names.Add('conversion_rate', '=Inputs!$B$10')

# --- 4) Add the "__names__" worksheet documenting all defined names --------
# Duplicate the README sheet (so the new sheet inherits clean default
# formatting/metadata) then clear it out and rename it.
$wsReadme = $wb.Worksheets.Item("README")
$wsReadme.Copy($null, $wsReadme)
$wsNames = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNames.Name = "__names__"
$wsNames.Cells.Clear()

$wsNames.Range("A1").Value = "Name"
$wsNames.Range("B1").Value = "RefersTo"

$wsNames.Range("A2").Value = "gross_margin"
$wsNames.Range("B2").Value = "Inputs!`$B`$2"

$wsNames.Range("A3").Value = "discount_rate"
$wsNames.Range("B3").Value = "Inputs!`$B`$3"

$wsNames.Range("A4").Value = "price_basic"
$wsNames.Range("B4").Value = "Inputs!`$B`$4"

$wsNames.Range("A5").Value = "price_pro"
$wsNames.Range("B5").Value = "Inputs!`$B`$5"

$wsNames.Range("A6").Value = "monthly_churn_basic"
$wsNames.Range("B6").Value = "Inputs!`$B`$6"

$wsNames.Range("A7").Value = "monthly_churn_pro"
$wsNames.Range("B7").Value = "Inputs!`$B`$7"

$wsNames.Range("A8").Value = "cac_target"
$wsNames.Range("B8").Value = "Inputs!`$B`$8"

$wsNames.Range("A9").Value = "activation_rate"
$wsNames.Range("B9").Value = "Inputs!`$B`$9"

$wsNames.Range("A10").Value = "conversion_rate"
$wsNames.Range("B10").Value = "Inputs!`$B`$10"

# Give the new header row the same bold/border/centered style used by the
# other sheets' header rows (copy formatting from the Inputs header cell).
$wsInputs.Range("A1").Copy()
$wsNames.Range("A1:B1").PasteSpecial(-4122)

# Keep the Inputs sheet as the active/selected tab (it was before our edits,
# and creating/activating the new sheet would otherwise steal that focus).
$wsInputs.Activate()
